$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 396 (existing data at row 396 and
# below shifts down by one, so old row 396 becomes row 397, etc.).
$ws.Rows.Item(396).Insert()

# Populate the newly inserted row 396 with the new weekly record
# (same market/product metadata as the row that used to sit there,
# but a new date and a new volume).
$ws.Cells.Item(396, 1).Value = 7
$ws.Cells.Item(396, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(396, 3).Value = "Ñuble"
$ws.Cells.Item(396, 4).Value = 44585
$ws.Cells.Item(396, 5).Value = 16
$ws.Cells.Item(396, 6).Value = "Fruta"
$ws.Cells.Item(396, 7).Value = 100106
$ws.Cells.Item(396, 8).Value = "Oleaginosos"
$ws.Cells.Item(396, 9).Value = 100106002
$ws.Cells.Item(396, 10).Value = "Palta"
$ws.Cells.Item(396, 11).Value = "Hass"
$ws.Cells.Item(396, 12).Value = "Primera"
$ws.Cells.Item(396, 13).Value = 120
$ws.Cells.Item(396, 14).Value = 2700
$ws.Cells.Item(396, 15).Value = 2800
$ws.Cells.Item(396, 16).Value = 2750
$ws.Cells.Item(396, 17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(396, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(396, 19).Value = 2750
$ws.Cells.Item(396, 20).Value = 1
